$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-8: existing rows get updated values (column B "Date" is unchanged, left untouched)
# Row 2
$ws.Range('A2').Value = 'Turkish 2 Lig'
$ws.Range('C2').Value = '07:00:00'
$ws.Range('D2').Value = 'Ankara Demirspor'
$ws.Range('E2').Value = 'Kahramanmaras Istiklal Spor'
$ws.Range('F2').Value = 1.8
$ws.Range('G2').Value = 1000
$ws.Range('H2').Value = 1.02
$ws.Range('I2').Value = 2.24
$ws.Range('J2').Value = 1.83
$ws.Range('K2').Value = 1000
$ws.Range('L2').Value = 1.28
$ws.Range('M2').Value = 1.01
$ws.Range('N2').Value = 1.81
$ws.Range('O2').Value = 1.01
$ws.Range('P2').Value = 1.81
$ws.Range('Q2').Value = 1.68
$ws.Range('R2').Value = 1.19
$ws.Range('S2').Value = 1.68
$ws.Range('T2').Value = 1.01
$ws.Range('U2').Value = 1.01
$ws.Range('V2').Value = 1.01
$ws.Range('W2').Value = 1.01
$ws.Range('X2').Value = 1000
$ws.Range('Y2').Value = 1000
$ws.Range('Z2').Value = 1000
$ws.Range('AA2').Value = 1000
$ws.Range('AB2').Value = 1000
$ws.Range('AC2').Value = 1000
$ws.Range('AD2').Value = 1000
$ws.Range('AE2').Value = 1000
$ws.Range('AF2').Value = 1000
$ws.Range('AG2').Value = 1000
$ws.Range('AH2').Value = 1000
$ws.Range('AI2').Value = 1000
$ws.Range('AJ2').Value = 1000
$ws.Range('AK2').Value = 1000
$ws.Range('AL2').Value = 1000
$ws.Range('AM2').Value = 1000
$ws.Range('AN2').Value = 1000
$ws.Range('AO2').Value = 1000

# Row 3
$ws.Range('A3').Value = 'Turkish 2 Lig'
$ws.Range('C3').Value = '07:00:00'
$ws.Range('D3').Value = 'Erbaaspor'
$ws.Range('E3').Value = 'Elazigspor'
$ws.Range('F3').Value = 1.02
$ws.Range('G3').Value = 1000
$ws.Range('H3').Value = 1.02
$ws.Range('I3').Value = 1000
$ws.Range('J3').Value = 1.02
$ws.Range('K3').Value = 1000
$ws.Range('L3').Value = 1.01
$ws.Range('M3').Value = 1.01
$ws.Range('N3').Value = 1.25
$ws.Range('O3').Value = 1.13
$ws.Range('P3').Value = 1.25
$ws.Range('Q3').Value = 1.13
$ws.Range('R3').Value = 1.18
$ws.Range('S3').Value = 1.17
$ws.Range('T3').Value = 1.01
$ws.Range('U3').Value = 1.01
$ws.Range('V3').Value = 1.01
$ws.Range('W3').Value = 1.01
$ws.Range('X3').Value = 1000
$ws.Range('Y3').Value = 1000
$ws.Range('Z3').Value = 1000
$ws.Range('AA3').Value = 1000
$ws.Range('AB3').Value = 1000
$ws.Range('AC3').Value = 1000
$ws.Range('AD3').Value = 1000
$ws.Range('AE3').Value = 1000
$ws.Range('AF3').Value = 1000
$ws.Range('AG3').Value = 1000
$ws.Range('AH3').Value = 1000
$ws.Range('AI3').Value = 1000
$ws.Range('AJ3').Value = 1000
$ws.Range('AK3').Value = 1000
$ws.Range('AL3').Value = 1000
$ws.Range('AM3').Value = 1000
$ws.Range('AN3').Value = 1000
$ws.Range('AO3').Value = 1000

# Row 4
$ws.Range('A4').Value = 'Turkish 2 Lig'
$ws.Range('C4').Value = '08:00:00'
$ws.Range('D4').Value = 'Mersin Idman Yurdu'
$ws.Range('E4').Value = 'Isparta 32 Spor'
$ws.Range('F4').Value = 3.6
$ws.Range('G4').Value = 6.8
$ws.Range('H4').Value = 1.69
$ws.Range('I4').Value = 2.22
$ws.Range('J4').Value = 3
$ws.Range('K4').Value = 7.6
$ws.Range('L4').Value = 1.33
$ws.Range('M4').Value = 1.01
$ws.Range('N4').Value = 1.66
$ws.Range('O4').Value = 1.01
$ws.Range('P4').Value = 1.66
$ws.Range('Q4').Value = 1.84
$ws.Range('R4').Value = 1.18
$ws.Range('S4').Value = 1.34
$ws.Range('T4').Value = 1.01
$ws.Range('U4').Value = 1.01
$ws.Range('V4').Value = 1.01
$ws.Range('W4').Value = 1.01
$ws.Range('X4').Value = 1000
$ws.Range('Y4').Value = 1000
$ws.Range('Z4').Value = 1000
$ws.Range('AA4').Value = 1000
$ws.Range('AB4').Value = 1000
$ws.Range('AC4').Value = 1000
$ws.Range('AD4').Value = 1000
$ws.Range('AE4').Value = 1000
$ws.Range('AF4').Value = 1000
$ws.Range('AG4').Value = 1000
$ws.Range('AH4').Value = 1000
$ws.Range('AI4').Value = 1000
$ws.Range('AJ4').Value = 1000
$ws.Range('AK4').Value = 1000
$ws.Range('AL4').Value = 1000
$ws.Range('AM4').Value = 1000
$ws.Range('AN4').Value = 1000
$ws.Range('AO4').Value = 1000

# Row 5
$ws.Range('A5').Value = 'Turkish 2 Lig'
$ws.Range('C5').Value = '08:00:00'
$ws.Range('D5').Value = 'Iskenderunspor'
$ws.Range('E5').Value = 'Adana 1954 FK'
$ws.Range('F5').Value = 1.02
$ws.Range('G5').Value = 1000
$ws.Range('H5').Value = 1.02
$ws.Range('I5').Value = 1000
$ws.Range('J5').Value = 1.02
$ws.Range('K5').Value = 1000
$ws.Range('L5').Value = 1.01
$ws.Range('M5').Value = 1.01
$ws.Range('N5').Value = 1.34
$ws.Range('O5').Value = 1.01
$ws.Range('P5').Value = 1.34
$ws.Range('Q5').Value = 1.01
$ws.Range('R5').Value = 1.18
$ws.Range('S5').Value = 1.37
$ws.Range('T5').Value = 1.01
$ws.Range('U5').Value = 1.01
$ws.Range('V5').Value = 1.01
$ws.Range('W5').Value = 1.01
$ws.Range('X5').Value = 1000
$ws.Range('Y5').Value = 1000
$ws.Range('Z5').Value = 1000
$ws.Range('AA5').Value = 1000
$ws.Range('AB5').Value = 1000
$ws.Range('AC5').Value = 1000
$ws.Range('AD5').Value = 1000
$ws.Range('AE5').Value = 1000
$ws.Range('AF5').Value = 1000
$ws.Range('AG5').Value = 1000
$ws.Range('AH5').Value = 1000
$ws.Range('AI5').Value = 1000
$ws.Range('AJ5').Value = 1000
$ws.Range('AK5').Value = 1000
$ws.Range('AL5').Value = 1000
$ws.Range('AM5').Value = 1000
$ws.Range('AN5').Value = 1000
$ws.Range('AO5').Value = 1000

# Row 6
$ws.Range('A6').Value = 'Turkish 2 Lig'
$ws.Range('C6').Value = '09:00:00'
$ws.Range('D6').Value = 'Bucaspor'
$ws.Range('E6').Value = 'Altinordu'
$ws.Range('F6').Value = 1.72
$ws.Range('G6').Value = 2.38
$ws.Range('H6').Value = 3.15
$ws.Range('I6').Value = 1000
$ws.Range('J6').Value = 2.96
$ws.Range('K6').Value = 1000
$ws.Range('L6').Value = 1.29
$ws.Range('M6').Value = 1.01
$ws.Range('N6').Value = 1.71
$ws.Range('O6').Value = 1.01
$ws.Range('P6').Value = 1.71
$ws.Range('Q6').Value = 1.74
$ws.Range('R6').Value = 1.18
$ws.Range('S6').Value = 1.74
$ws.Range('T6').Value = 1.01
$ws.Range('U6').Value = 1.01
$ws.Range('V6').Value = 1.32
$ws.Range('W6').Value = 1.72
$ws.Range('X6').Value = 1000
$ws.Range('Y6').Value = 1000
$ws.Range('Z6').Value = 1000
$ws.Range('AA6').Value = 1000
$ws.Range('AB6').Value = 1000
$ws.Range('AC6').Value = 1000
$ws.Range('AD6').Value = 1000
$ws.Range('AE6').Value = 1000
$ws.Range('AF6').Value = 1000
$ws.Range('AG6').Value = 1000
$ws.Range('AH6').Value = 1000
$ws.Range('AI6').Value = 1000
$ws.Range('AJ6').Value = 1000
$ws.Range('AK6').Value = 1000
$ws.Range('AL6').Value = 1000
$ws.Range('AM6').Value = 1000
$ws.Range('AN6').Value = 1000
$ws.Range('AO6').Value = 1000

# Row 7
$ws.Range('A7').Value = 'Turkish 2 Lig'
$ws.Range('C7').Value = '09:00:00'
$ws.Range('D7').Value = 'Beykoz Anadolu Spor'
$ws.Range('E7').Value = 'Batman Petrolspor'
$ws.Range('F7').Value = 4.8
$ws.Range('G7').Value = 1000
$ws.Range('H7').Value = 1.45
$ws.Range('I7').Value = 1.85
$ws.Range('J7').Value = 3.35
$ws.Range('K7').Value = 9.6
$ws.Range('L7').Value = 1.01
$ws.Range('M7').Value = 1.01
$ws.Range('N7').Value = 1.75
$ws.Range('O7').Value = 1.01
$ws.Range('P7').Value = 1.75
$ws.Range('Q7').Value = 1.74
$ws.Range('R7').Value = 1.19
$ws.Range('S7').Value = 1.74
$ws.Range('T7').Value = 1.01
$ws.Range('U7').Value = 1.01
$ws.Range('V7').Value = 2.16
$ws.Range('W7').Value = 1.01
$ws.Range('X7').Value = 1000
$ws.Range('Y7').Value = 1000
$ws.Range('Z7').Value = 1000
$ws.Range('AA7').Value = 1000
$ws.Range('AB7').Value = 1000
$ws.Range('AC7').Value = 1000
$ws.Range('AD7').Value = 1000
$ws.Range('AE7').Value = 1000
$ws.Range('AF7').Value = 1000
$ws.Range('AG7').Value = 1000
$ws.Range('AH7').Value = 1000
$ws.Range('AI7').Value = 1000
$ws.Range('AJ7').Value = 1000
$ws.Range('AK7').Value = 1000
$ws.Range('AL7').Value = 1000
$ws.Range('AM7').Value = 1000
$ws.Range('AN7').Value = 1000
$ws.Range('AO7').Value = 1000

# Row 8
$ws.Range('A8').Value = 'Serbian Super League'
$ws.Range('C8').Value = '12:30:00'
$ws.Range('D8').Value = 'FK Radnicki 1923'
$ws.Range('E8').Value = 'Cukaricki'
$ws.Range('F8').Value = 2.14
$ws.Range('G8').Value = 2.74
$ws.Range('H8').Value = 2.88
$ws.Range('I8').Value = 3.85
$ws.Range('J8').Value = 3.15
$ws.Range('K8').Value = 4.3
$ws.Range('L8').Value = 1.27
$ws.Range('M8').Value = 1.05
$ws.Range('N8').Value = 3.75
$ws.Range('O8').Value = 1.23
$ws.Range('P8').Value = 2.06
$ws.Range('Q8').Value = 1.62
$ws.Range('R8').Value = 1.41
$ws.Range('S8').Value = 2.56
$ws.Range('T8').Value = 1.59
$ws.Range('U8').Value = 2.22
$ws.Range('V8').Value = 1.37
$ws.Range('W8').Value = 1.57
$ws.Range('X8').Value = 1000
$ws.Range('Y8').Value = 1000
$ws.Range('Z8').Value = 1000
$ws.Range('AA8').Value = 1000
$ws.Range('AB8').Value = 1000
$ws.Range('AC8').Value = 1000
$ws.Range('AD8').Value = 1000
$ws.Range('AE8').Value = 1000
$ws.Range('AF8').Value = 1000
$ws.Range('AG8').Value = 1000
$ws.Range('AH8').Value = 1000
$ws.Range('AI8').Value = 1000
$ws.Range('AJ8').Value = 1000
$ws.Range('AK8').Value = 1000
$ws.Range('AL8').Value = 1000
$ws.Range('AM8').Value = 1000
$ws.Range('AN8').Value = 1000
$ws.Range('AO8').Value = 1000

# Rows 9-12: brand-new rows appended at the end of the data block
# Row 9
$ws.Range('A9').Value = 'Swiss Super League'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = '2025-12-17'
$ws.Range('C9').Value = '16:30:00'
$ws.Range('D9').Value = 'Young Boys'
$ws.Range('E9').Value = 'Grasshoppers Zurich'
$ws.Range('F9').Value = 1.61
$ws.Range('G9').Value = 1.64
$ws.Range('H9').Value = 5.4
$ws.Range('I9').Value = 6
$ws.Range('J9').Value = 4.7
$ws.Range('K9').Value = 5.1
$ws.Range('L9').Value = 1.01
$ws.Range('M9').Value = 1.03
$ws.Range('N9').Value = 6.2
$ws.Range('O9').Value = 1.16
$ws.Range('P9').Value = 2.8
$ws.Range('Q9').Value = 1.49
$ws.Range('R9').Value = 1.74
$ws.Range('S9').Value = 2.24
$ws.Range('T9').Value = 1.61
$ws.Range('U9').Value = 2.44
$ws.Range('V9').Value = 1.2
$ws.Range('W9').Value = 2.56
$ws.Range('X9').Value = 29
$ws.Range('Y9').Value = 980
$ws.Range('Z9').Value = 55
$ws.Range('AA9').Value = 140
$ws.Range('AB9').Value = 15
$ws.Range('AC9').Value = 12
$ws.Range('AD9').Value = 980
$ws.Range('AE9').Value = 60
$ws.Range('AF9').Value = 13
$ws.Range('AG9').Value = 10.5
$ws.Range('AH9').Value = 17.5
$ws.Range('AI9').Value = 55
$ws.Range('AJ9').Value = 17.5
$ws.Range('AK9').Value = 15
$ws.Range('AL9').Value = 25
$ws.Range('AM9').Value = 70
$ws.Range('AN9').Value = 5.7
$ws.Range('AO9').Value = 44

# Row 10
$ws.Range('A10').Value = 'Swiss Super League'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = '2025-12-17'
$ws.Range('C10').Value = '16:30:00'
$ws.Range('D10').Value = 'FC Zurich'
$ws.Range('E10').Value = 'Lugano'
$ws.Range('F10').Value = 2.74
$ws.Range('G10').Value = 3.05
$ws.Range('H10').Value = 2.46
$ws.Range('I10').Value = 2.68
$ws.Range('J10').Value = 3.6
$ws.Range('K10').Value = 4
$ws.Range('L10').Value = 1.29
$ws.Range('M10').Value = 1.05
$ws.Range('N10').Value = 4.4
$ws.Range('O10').Value = 1.25
$ws.Range('P10').Value = 2.16
$ws.Range('Q10').Value = 1.73
$ws.Range('R10').Value = 1.46
$ws.Range('S10').Value = 2.58
$ws.Range('T10').Value = 1.63
$ws.Range('U10').Value = 2.34
$ws.Range('V10').Value = 1.6
$ws.Range('W10').Value = 1.5
$ws.Range('X10').Value = 19
$ws.Range('Y10').Value = 13.5
$ws.Range('Z10').Value = 19.5
$ws.Range('AA10').Value = 1000
$ws.Range('AB10').Value = 14.5
$ws.Range('AC10').Value = 9
$ws.Range('AD10').Value = 13
$ws.Range('AE10').Value = 28
$ws.Range('AF10').Value = 22
$ws.Range('AG10').Value = 14
$ws.Range('AH10').Value = 16.5
$ws.Range('AI10').Value = 1000
$ws.Range('AJ10').Value = 1000
$ws.Range('AK10').Value = 1000
$ws.Range('AL10').Value = 1000
$ws.Range('AM10').Value = 1000
$ws.Range('AN10').Value = 24
$ws.Range('AO10').Value = 19

# Row 11
$ws.Range('A11').Value = 'Swiss Super League'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = '2025-12-17'
$ws.Range('C11').Value = '16:30:00'
$ws.Range('D11').Value = 'Luzern'
$ws.Range('E11').Value = 'FC Basel'
$ws.Range('F11').Value = 3.45
$ws.Range('G11').Value = 3.65
$ws.Range('H11').Value = 2.08
$ws.Range('I11').Value = 2.18
$ws.Range('J11').Value = 3.85
$ws.Range('K11').Value = 4.2
$ws.Range('L11').Value = 1.25
$ws.Range('M11').Value = 1.03
$ws.Range('N11').Value = 5.6
$ws.Range('O11').Value = 1.18
$ws.Range('P11').Value = 2.52
$ws.Range('Q11').Value = 1.55
$ws.Range('R11').Value = 1.64
$ws.Range('S11').Value = 2.3
$ws.Range('T11').Value = 1.53
$ws.Range('U11').Value = 2.66
$ws.Range('V11').Value = 1.84
$ws.Range('W11').Value = 1.37
$ws.Range('X11').Value = 1000
$ws.Range('Y11').Value = 1000
$ws.Range('Z11').Value = 1000
$ws.Range('AA11').Value = 980
$ws.Range('AB11').Value = 1000
$ws.Range('AC11').Value = 10
$ws.Range('AD11').Value = 12
$ws.Range('AE11').Value = 1000
$ws.Range('AF11').Value = 1000
$ws.Range('AG11').Value = 16
$ws.Range('AH11').Value = 15
$ws.Range('AI11').Value = 1000
$ws.Range('AJ11').Value = 60
$ws.Range('AK11').Value = 980
$ws.Range('AL11').Value = 1000
$ws.Range('AM11').Value = 55
$ws.Range('AN11').Value = 1000
$ws.Range('AO11').Value = 9.4

# Row 12
$ws.Range('A12').Value = 'Scottish Premiership'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = '2025-12-17'
$ws.Range('C12').Value = '17:00:00'
$ws.Range('D12').Value = 'Dundee Utd'
$ws.Range('E12').Value = 'Celtic'
$ws.Range('F12').Value = 7.4
$ws.Range('G12').Value = 7.8
$ws.Range('H12').Value = 1.5
$ws.Range('I12').Value = 1.52
$ws.Range('J12').Value = 4.8
$ws.Range('K12').Value = 5
$ws.Range('L12').Value = 1.31
$ws.Range('M12').Value = 1.04
$ws.Range('N12').Value = 5.2
$ws.Range('O12').Value = 1.21
$ws.Range('P12').Value = 2.42
$ws.Range('Q12').Value = 1.61
$ws.Range('R12').Value = 1.59
$ws.Range('S12').Value = 2.54
$ws.Range('T12').Value = 1.77
$ws.Range('U12').Value = 2.14
$ws.Range('V12').Value = 2.92
$ws.Range('W12').Value = 1.14
$ws.Range('X12').Value = 23
$ws.Range('Y12').Value = 11
$ws.Range('Z12').Value = 10.5
$ws.Range('AA12').Value = 14
$ws.Range('AB12').Value = 32
$ws.Range('AC12').Value = 11.5
$ws.Range('AD12').Value = 10.5
$ws.Range('AE12').Value = 14.5
$ws.Range('AF12').Value = 60
$ws.Range('AG12').Value = 26
$ws.Range('AH12').Value = 21
$ws.Range('AI12').Value = 30
$ws.Range('AJ12').Value = 220
$ws.Range('AK12').Value = 90
$ws.Range('AL12').Value = 85
$ws.Range('AM12').Value = 110
$ws.Range('AN12').Value = 90
$ws.Range('AO12').Value = 6.4
